# Update the "Förändrad" (Changed) date column (C) for rows 2-11
# from 2023-09-01 (45170) to 2023-09-05 (45174), preserving existing
# cell formatting/style.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 11; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45170) {
        $cell.Value = 45174
    }
}
